$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Datos actualizados a 21 de Marzo de 2020 a las 05:46"

$ws.Range("B8").Value = 19658
$ws.Range("C8").Value = 275
$ws.Range("E8").Value = 19247

$ws.Range("A55").Value = "Mexico"
$ws.Range("B55").Value = 203
$ws.Range("C55").Value = 39
$ws.Range("D55").Value = 4
$ws.Range("E55").Value = 197
$ws.Range("F55").Value = 1
$ws.Range("G55").Value = 1
$ws.Range("H55").Value = 2

$ws.Range("A56").Value = "Sudafrica"
$ws.Range("B56").Value = 202
$ws.Range("C56").Value = 0
$ws.Range("D56").Value = 0
$ws.Range("E56").Value = 202
$ws.Range("F56").Value = 0
$ws.Range("G56").Value = 0
$ws.Range("H56").Value = 0

$ws.Range("A57").Value = "Panama"
$ws.Range("B57").Value = 200
$ws.Range("C57").Value = 0
$ws.Range("D57").Value = 1
$ws.Range("E57").Value = 198
$ws.Range("F57").Value = 7
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 1

$ws.Range("A58").Value = "Libano"
$ws.Range("B58").Value = 177
$ws.Range("C58").Value = 0
$ws.Range("D58").Value = 4
$ws.Range("E58").Value = 169
$ws.Range("F58").Value = 3
$ws.Range("G58").Value = 0
$ws.Range("H58").Value = 4

$ws.Range("A144").Value = "Islas Virgenes de los Estados Unidos"
$ws.Range("B144").Value = 6
$ws.Range("C144").Value = 3
$ws.Range("D144").Value = 0
$ws.Range("E144").Value = 6
$ws.Range("F144").Value = 0
$ws.Range("G144").Value = 0
$ws.Range("H144").Value = 0

$ws.Range("A145").Value = "Aruba"
$ws.Range("B145").Value = 5
$ws.Range("C145").Value = 0
$ws.Range("D145").Value = 1
$ws.Range("E145").Value = 4
$ws.Range("F145").Value = 0
$ws.Range("G145").Value = 0
$ws.Range("H145").Value = 0

$ws.Range("A148").Value = "Bahamas"
$ws.Range("B148").Value = 4
$ws.Range("C148").Value = 0
$ws.Range("D148").Value = 0
$ws.Range("E148").Value = 4
$ws.Range("F148").Value = 0
$ws.Range("G148").Value = 0
$ws.Range("H148").Value = 0

$ws.Range("A149").Value = "Gabon"
$ws.Range("B149").Value = 4
$ws.Range("C149").Value = 0
$ws.Range("D149").Value = 0
$ws.Range("E149").Value = 3
$ws.Range("F149").Value = 0
$ws.Range("G149").Value = 0
$ws.Range("H149").Value = 1

$ws.Range("A150").Value = "Madagascar"

$ws.Range("A152").Value = "Republica de Africa Central"
$ws.Range("A153").Value = "San Bartolome"
$ws.Range("A154").Value = "Congo"

$ws.Range("A155").Value = "Islas Caimanes"
$ws.Range("A156").Value = "Curazao"

$ws.Range("A157").Value = "Liberia"
$ws.Range("A158").Value = "Santa Lucia"
$ws.Range("A159").Value = "Zambia"
$ws.Range("A160").Value = "Haiti"
$ws.Range("A161").Value = "Butan"
$ws.Range("A162").Value = "Nueva Caledonia"
$ws.Range("A163").Value = "Benin"
$ws.Range("A164").Value = "Bermudas"
$ws.Range("A165").Value = "Isla de Man"
$ws.Range("A166").Value = "Nicaragua"
$ws.Range("A167").Value = "Guinea"
$ws.Range("A168").Value = "Groenlandia"
$ws.Range("A169").Value = "Mauritania"
